$d = $word.ActiveDocument
$d.Content.Find.Execute("Claude Code,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Claude Code, ", 2)
